$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.410.28'

# Row 3
$ws.Range('D3').Value = '1.723.73'
$ws.Range('E3').Value = '  -0.38%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.22'
$ws.Range('E5').Value = '  -1.19%  '

# Row 6
$ws.Range('E6').Value = '  +0.04%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4911'
$ws.Range('E7').Value = '  +1.86%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2605'
$ws.Range('E8').Value = '  -2.50%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06198'
$ws.Range('E9').Value = '  +0.26%  '

# Row 10
$ws.Range('D10').Value = '1.717.03'
$ws.Range('E10').Value = '  -0.97%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06988'
$ws.Range('E11').Value = '  -1.61%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.44'
$ws.Range('E12').Value = '  -1.29%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.543'
$ws.Range('E13').Value = '  +0.00%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5991'
$ws.Range('E14').Value = '  -2.32%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.38'
$ws.Range('E15').Value = '  +0.09%  '

# Row 16
$ws.Range('E16').Value = '  +0.02%  '

# Row 17
$ws.Range('D17').Value = '26.406.81'

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9998'
$ws.Range('E18').Value = '  +0.00%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007215'
$ws.Range('E19').Value = '  +3.63%  '

# Row 21
$ws.Range('D21').Value = '1.944.20'
$ws.Range('E21').Value = '  -0.57%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.465'
$ws.Range('E22').Value = '  -1.35%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.582'
$ws.Range('E23').Value = '  -2.71%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.154'
$ws.Range('E24').Value = '  -1.78%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '137.52'
$ws.Range('E25').Value = '  +0.03%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.400'
$ws.Range('E27').Value = '  -0.45%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '106.90'
$ws.Range('E28').Value = '  -1.20%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.721'
$ws.Range('E29').Value = '  -3.29%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.947'
$ws.Range('E30').Value = '  -0.87%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07999'
$ws.Range('E31').Value = '  -0.25%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.684'
$ws.Range('E32').Value = '  -0.18%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04505'
$ws.Range('E33').Value = '  -0.87%  '

# Row 34
$ws.Range('B34').Value = 'Frax'
$ws.Range('C34').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9993'
$ws.Range('E34').Value = '  +0.01%  '

# Row 35
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.603'
$ws.Range('E35').Value = '  -0.41%  '

# Row 36
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9998'
$ws.Range('E36').Value = '  -0.61%  '

# Row 37
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6258'
$ws.Range('E37').Value = '  -1.13%  '

# Row 38
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9413'
$ws.Range('E38').Value = '  +4.70%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.390'
$ws.Range('E39').Value = '  +0.02%  '

# Row 40
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.945'
$ws.Range('E40').Value = '  -5.13%  '

# Row 41
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9998'
$ws.Range('E41').Value = '  -0.39%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.01481'
$ws.Range('E42').Value = '  -1.40%  '

# Row 43
$ws.Range('B43').Value = 'Quant'
$ws.Range('C43').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '99.65'
$ws.Range('E43').Value = '  -3.16%  '

# Row 44
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.307'
$ws.Range('E44').Value = '  -2.53%  '

# Row 45
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.3850'
$ws.Range('E45').Value = '  -1.40%  '

# Row 46
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.816'
$ws.Range('E46').Value = '  -5.09%  '

# Row 47
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.1169'
$ws.Range('E47').Value = '  -1.60%  '

# Row 48
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05362'
$ws.Range('E48').Value = '  -0.51%  '

# Row 49
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.732'
$ws.Range('E49').Value = '  -2.20%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '30.17'
$ws.Range('E50').Value = '  -1.64%  '

# Row 51
$ws.Range('B51').Value = 'NEARProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.235'
$ws.Range('E51').Value = '  -1.71%  '
